# "Working on next button" - update rental property listings to the next page of results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep bedrooms/bathrooms columns as text (matching original inlineStr content,
# e.g. "3", "4 + 1") instead of letting Excel auto-coerce them to numbers.
$ws.Range("C2:D13").NumberFormat = "@"

$data = @(
    @("232 NICKERSON DRIVE, Cobourg, Ontario", "$2,700/Monthly", "3", "2"),
    @("LOT 114 - 1081 DENTON DRIVE, Cobourg, Ontario", "$3,500/Monthly", "4 + 1", "4"),
    @("404 - 325 UNIVERSITY AVENUE W, Cobourg, Ontario", "$2,025/Monthly", "1", "1"),
    @("1016 TRAILSVIEW AVENUE, Cobourg, Ontario", "$3,499/Monthly", "4", "4"),
    @("UNIT 1 - 74 KING STREET W, Cobourg, Ontario", "$1,600/Monthly", "1", "1"),
    @("319 - 325 UNIVERSITY AVENUE W, Cobourg, Ontario", "$2,050/Monthly", "1", "1"),
    @("1004 TRAILSVIEW AVENUE, Cobourg, Ontario", "$3,500/Monthly", "4", "5"),
    @("417 TREVOR STREET, Cobourg, Ontario", "$2,900/Monthly", "3 + 1", "3"),
    @("502 - 79 KING STREET, Cobourg, Ontario", "$2,650/Monthly", "3", "2"),
    @("475 DREWERY ROAD, Cobourg, Ontario", "$3,200/Monthly", "3", "3"),
    @("27 - 160 DENSMORE ROAD, Cobourg, Ontario", "$2,890/Monthly", "3", "3"),
    @("1066 DENTON DRIVE, Cobourg, Ontario", "$3,400/Monthly", "5", "4")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
